$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45996, 10691.1282886068, 9304.0772333742, 16128.26, 6366.53634899879, -19.0686007344588),
    @(45997, 4139.00616187948, 6693.95511100164, 11232.26, 6439.9709402312, 79.2360854680353),
    @(45998, 3964.17494499607, 6619.05729364572, 11232.26, 6409.92321032582, 74.8633543321478),
    @(45999, 4099.92020470063, 6558.73357681079, 11232.26, 6651.94893001777, 82.4342711178567),
    @(46000, 10890.0129278173, 10211.269508972, 11232.26, 6750.54030538153, 238.731242264729),
    @(46001, 9706.9828974431, 9745.24918392598, 11232.26, 7769.40825216837, 261.766559837265),
    @(46002, 9706.9828974431, 9762.0588145387, 11232.26, 7769.40825216837, 262.466961112794),
    @(46003, 9706.9828974431, 8977.00331093427, 11232.26, 7769.40825216837, 229.756315129276),
    @(46004, 8488.30621711765, 8504.08867310034, 11232.26, 7366.61649868131, 193.268548824235),
    @(46005, 8386.73995792576, 8276.61231181288, 11232.26, 7357.49801940159, 183.41043046727),
    @(46006, 9746.03784374469, 8944.32407478398, 11232.26, 8042.20026808171, 239.76101428607),
    @(46007, 9746.03784374469, 9241.70378611331, 11232.26, 8042.20026808171, 252.151835591459),
    @(46008, 9746.03784374469, 9270.65930142044, 11232.26, 8042.20026808171, 253.358315395923),
    @(46009, 9746.03784374469, 9260.20653093353, 11232.26, 8042.20026808171, 252.922783292302)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $rowIndex++
}
